{"js": "// Lattice multiplication exercises: replace the 15 problems (5 rows x 3 cols)\n// in the single table with a new set of multiplicand pairs, keeping each\n// cell's existing \"lattice\" layout/formatting:\n//   AB x CD\n//     C    D\n//   ----\n//   A|    |\n//   B|    |\n\nconst pairs = [\n  [2, 9, 8, 5], [2, 9, 6, 6], [9, 0, 4, 4],\n  [3, 9, 8, 7], [7, 8, 7, 9], [8, 6, 4, 0],\n  [4, 3, 3, 0], [7, 9, 5, 0], [7, 7, 4, 3],\n  [2, 9, 8, 9], [7, 5, 7, 9], [5, 9, 3, 3],\n  [5, 3, 7, 4], [3, 5, 1, 0], [5, 3, 2, 5],\n];\n\nfunction cellText(a, b, c, d) {\n  return [\n    `${a}${b} x ${c}${d}`,\n    `  ${c}    ${d}`,\n    `  ----`,\n    `${a}|    |`,\n    `${b}|    |`,\n  ].join(\"\\u000b\"); // \\u000b (vertical tab) -> <w:br/> line break within the run\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst cols = 3;\nconst rows = pairs.length / cols;\n\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const [a, b, cc, d] = pairs[r * cols + c];\n    const cell = table.getCell(r, c);\n    // Replace the cell's whole range in place so the existing run\n    // formatting (e.g. sz=32) on that paragraph is kept.\n    const range = cell.body.getRange();\n    range.insertText(cellText(a, b, cc, d), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication exercises: replace the 15 problems (5 rows x 3 cols)\n# in the single table with a new set of multiplicand pairs, keeping each\n# cell's existing \"lattice\" layout/formatting:\n#   AB x CD\n#     C    D\n#   ----\n#   A|    |\n#   B|    |\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairs = @(\n    @(2, 9, 8, 5), @(2, 9, 6, 6), @(9, 0, 4, 4),\n    @(3, 9, 8, 7), @(7, 8, 7, 9), @(8, 6, 4, 0),\n    @(4, 3, 3, 0), @(7, 9, 5, 0), @(7, 7, 4, 3),\n    @(2, 9, 8, 9), @(7, 5, 7, 9), @(5, 9, 3, 3),\n    @(5, 3, 7, 4), @(3, 5, 1, 0), @(5, 3, 2, 5)\n)\n\n$lineBreak = [char]11   # -> <w:br/> inside the run, same story as a Word soft line break\n$cols = 3\n\nfor ($i = 0; $i -lt $pairs.Count; $i++) {\n    $pair = $pairs[$i]\n    $a = $pair[0]\n    $b = $pair[1]\n    $c = $pair[2]\n    $dd = $pair[3]\n    $row = [int][math]::Floor($i / $cols) + 1\n    $col = ($i % $cols) + 1\n\n    $lines = @(\n        \"$a$b x $c$dd\",\n        \"  $c    $dd\",\n        \"  ----\",\n        \"$a|    |\",\n        \"$b|    |\"\n    )\n    $text = [string]::Join($lineBreak, $lines)\n\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
